# Sprint 44 test case report — fill in Day 8, Day 9 and Day 10 summary
# numbers (windows check for new release and ui test case written).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Spint( 44) - Day 8 - Test Case Summary (rows 39-41)
$ws.Range("C39").Value = 7075
$ws.Range("C40").Value = 2790
$ws.Range("C41").Value = 2790

# Spint( 44) - Day 9 - Test Case Summary (rows 45-47)
$ws.Range("C45").Value = 7103
$ws.Range("C46").Value = 2850
$ws.Range("C47").Value = 2850

# Spint( 44) - Day 10 - Test Case Summary (rows 51-53)
$ws.Range("C51").Value = 7131
$ws.Range("C52").Value = 2880
$ws.Range("C53").Value = 2880

# Scroll the view down to the newly-entered data and leave the last
# touched cell (C53) selected, matching where the author ended up.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C53").Select()
